# July 16 - email password update and run inputs update
# Updates the FBPixels sheet: change one campaign name, remove two obsolete
# rows, and make FBPixels the active/selected sheet (instead of AllPixels).

$wb = $excel.ActiveWorkbook

$wsFB = $wb.Worksheets.Item("FBPixels")
$wsAll = $wb.Worksheets.Item("AllPixels")

# Update the Sub-D / deluxe25offp row to the new July4 campaign name.
$wsFB.Range("C9").Value = "deluxe25offp-redes-July4"

# Remove the now-obsolete "special-offer" (row 14) and "deluxe25off" (row 12)
# rows - delete the higher row index first so the other index stays valid.
$wsFB.Rows.Item(14).Delete()
$wsFB.Rows.Item(12).Delete()

# Make FBPixels the active sheet/tab with a fresh selection, and drop the
# previous AllPixels active-tab state (its own selection is unchanged).
$wsAll.Range("C7").Select()
$wsFB.Activate()
$wsFB.Range("C11").Select()
